$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are stored as text, matching the source data
# (these look numeric but must preserve exact formatting, e.g. trailing zeros,
# thousand-dot separators, and subscript characters).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.653.38"
$ws.Range("E2").Value = "  +8.88%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.581.18"
$ws.Range("E3").Value = "  +10.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "503.31"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.99"
$ws.Range("E6").Value = "  +9.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("E7").Value = "  +22.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.580.77"
$ws.Range("E9").Value = "  +9.99%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.14"
$ws.Range("E10").Value = "  +13.48%  "

$ws.Range("E11").Value = "  +6.74%  "

$ws.Range("E12").Value = "  +6.58%  "

$ws.Range("E13").Value = "  +1.80%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.006.75"
$ws.Range("E14").Value = "  +9.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.508.00"
$ws.Range("E15").Value = "  +8.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.68"
$ws.Range("E16").Value = "  +8.70%  "

$ws.Range("E17").Value = "  +5.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.574.16"
$ws.Range("E18").Value = "  +9.44%  "

$ws.Range("E19").Value = "  +4.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "331.74"
$ws.Range("E20").Value = "  +6.43%  "

$ws.Range("E21").Value = "  +7.95%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.03"
$ws.Range("E22").Value = "  +7.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.01"
$ws.Range("E23").Value = "  +0.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "59.81"
$ws.Range("E24").Value = "  +6.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.414"
$ws.Range("E25").Value = "  +6.25%  "

$ws.Range("E26").Value = "  +8.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.993"
$ws.Range("E27").Value = "  -0.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.653.73"
$ws.Range("E28").Value = "  +8.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.39"
$ws.Range("E29").Value = "  +3.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0822"
$ws.Range("E30").Value = "  +10.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.38"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.28"
$ws.Range("E33").Value = "  +7.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.56"
$ws.Range("E34").Value = "  +6.97%  "

$ws.Range("E35").Value = "  +9.70%  "

$ws.Range("E36").Value = "  +9.66%  "

$ws.Range("E37").Value = "  +9.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.840"
$ws.Range("E38").Value = "  +3.03%  "

$ws.Range("E39").Value = "  +11.95%  "

$ws.Range("E40").Value = "  +8.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.11"
$ws.Range("E41").Value = "  +5.38%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "290.78"
$ws.Range("E42").Value = "  +16.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.101"
$ws.Range("E43").Value = "  +6.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.625"
$ws.Range("E44").Value = "  +8.80%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0562"
$ws.Range("E45").Value = "  +7.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("B46").Value = "SuiNetwork"
$ws.Range("C46").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D46").Value = "0.775"
$ws.Range("E46").Value = "  +24.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("E48").Value = "  +14.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0235"
$ws.Range("E49").Value = "  +7.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "4.79"
$ws.Range("E50").Value = "  +9.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.006.51"
$ws.Range("E51").Value = "  +12.64%  "
